$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 335 for a new week of price reports (2021-11-09),
# shifting the existing rows 335:349 down to 338:352
$ws.Rows.Item(335).Insert()
$ws.Rows.Item(335).Insert()
$ws.Rows.Item(335).Insert()

# Common constant values shared by every data row in this sheet
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112006
$categoria = "Repollo"
$unidad = "`$/unidad"
$kgUnidades = 1
$clasificacion = "Hortaliza"
$fecha = 44509

# New row 335
$r = 335
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Crespo record"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 5200
$ws.Cells.Item($r, 11).Value = 600
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 650
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 650
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 336
$r = 336
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Crespo record"
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 2500
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 337
$r = 337
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Morada(o)"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 2500
$ws.Cells.Item($r, 11).Value = 700
$ws.Cells.Item($r, 12).Value = 800
$ws.Cells.Item($r, 13).Value = 750
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 16).Value = 750
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

$ws.Range("A1").Select()
